# Auto-generated script to apply NBA team-data corrections
# (data values taken from a corrected 1-day date alignment + date format fix)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric/statistic corrections ---
$ws.Cells.Item(2, 30).Value = 3
$ws.Cells.Item(3, 30).Value = 14
$ws.Cells.Item(3, 37).Value = 9
$ws.Cells.Item(3, 40).Value = 18
$ws.Cells.Item(4, 30).Value = 14
$ws.Cells.Item(5, 30).Value = 14
$ws.Cells.Item(5, 34).Value = 17
$ws.Cells.Item(5, 45).Value = 28
$ws.Cells.Item(6, 30).Value = 26
$ws.Cells.Item(6, 37).Value = 25
$ws.Cells.Item(6, 42).Value = 23
$ws.Cells.Item(7, 30).Value = 26
$ws.Cells.Item(8, 4).Value = 71
$ws.Cells.Item(8, 6).Value = 36
$ws.Cells.Item(8, 7).Value = 0.493
$ws.Cells.Item(8, 9).Value = 38.7
$ws.Cells.Item(8, 11).Value = 0.461
$ws.Cells.Item(8, 14).Value = 0.374
$ws.Cells.Item(8, 15).Value = 16.8
$ws.Cells.Item(8, 16).Value = 21.1
$ws.Cells.Item(8, 17).Value = 0.795
$ws.Cells.Item(8, 19).Value = 32.7
$ws.Cells.Item(8, 20).Value = 42
$ws.Cells.Item(8, 21).Value = 23.1
$ws.Cells.Item(8, 28).Value = 101.8
$ws.Cells.Item(8, 29).Value = -0.5
$ws.Cells.Item(8, 30).Value = 14
$ws.Cells.Item(8, 32).Value = 16
$ws.Cells.Item(8, 33).Value = 17
$ws.Cells.Item(8, 37).Value = 8
$ws.Cells.Item(8, 39).Value = 11
$ws.Cells.Item(8, 42).Value = 21
$ws.Cells.Item(10, 30).Value = 3
$ws.Cells.Item(10, 40).Value = 17
$ws.Cells.Item(10, 50).Value = 15
$ws.Cells.Item(10, 55).Value = 26
$ws.Cells.Item(11, 34).Value = 21
$ws.Cells.Item(11, 46).Value = 4
$ws.Cells.Item(12, 30).Value = 14
$ws.Cells.Item(12, 35).Value = 6
$ws.Cells.Item(13, 4).Value = 72
$ws.Cells.Item(13, 5).Value = 45
$ws.Cells.Item(13, 7).Value = 0.625
$ws.Cells.Item(13, 10).Value = 80.90000000000001
$ws.Cells.Item(13, 11).Value = 0.434
$ws.Cells.Item(13, 13).Value = 19.7
$ws.Cells.Item(13, 14).Value = 0.355
$ws.Cells.Item(13, 15).Value = 17.3
$ws.Cells.Item(13, 16).Value = 22.9
$ws.Cells.Item(13, 17).Value = 0.752
$ws.Cells.Item(13, 18).Value = 12.8
$ws.Cells.Item(13, 20).Value = 46.3
$ws.Cells.Item(13, 21).Value = 20.4
$ws.Cells.Item(13, 23).Value = 7
$ws.Cells.Item(13, 27).Value = 21.6
$ws.Cells.Item(13, 28).Value = 94.5
$ws.Cells.Item(13, 29).Value = 4.9
$ws.Cells.Item(13, 30).Value = 3
$ws.Cells.Item(13, 33).Value = 8
$ws.Cells.Item(13, 34).Value = 12
$ws.Cells.Item(13, 37).Value = 28
$ws.Cells.Item(13, 40).Value = 20
$ws.Cells.Item(13, 48).Value = 23
$ws.Cells.Item(14, 30).Value = 3
$ws.Cells.Item(14, 40).Value = 19
$ws.Cells.Item(15, 4).Value = 72
$ws.Cells.Item(15, 6).Value = 35
$ws.Cells.Item(15, 7).Value = 0.514
$ws.Cells.Item(15, 9).Value = 37.3
$ws.Cells.Item(15, 13).Value = 24.6
$ws.Cells.Item(15, 14).Value = 0.36
$ws.Cells.Item(15, 15).Value = 18.9
$ws.Cells.Item(15, 16).Value = 27.5
$ws.Cells.Item(15, 17).Value = 0.6870000000000001
$ws.Cells.Item(15, 19).Value = 33.3
$ws.Cells.Item(15, 20).Value = 44.8
$ws.Cells.Item(15, 21).Value = 22.1
$ws.Cells.Item(15, 22).Value = 15.1
$ws.Cells.Item(15, 29).Value = 0.9
$ws.Cells.Item(15, 30).Value = 3
$ws.Cells.Item(15, 34).Value = 28
$ws.Cells.Item(15, 42).Value = 2
$ws.Cells.Item(15, 46).Value = 3
$ws.Cells.Item(16, 30).Value = 14
$ws.Cells.Item(16, 34).Value = 17
$ws.Cells.Item(16, 42).Value = 22
$ws.Cells.Item(17, 30).Value = 14
$ws.Cells.Item(17, 47).Value = 12
$ws.Cells.Item(18, 4).Value = 70
$ws.Cells.Item(18, 5).Value = 34
$ws.Cells.Item(18, 7).Value = 0.486
$ws.Cells.Item(18, 9).Value = 38.1
$ws.Cells.Item(18, 10).Value = 87.40000000000001
$ws.Cells.Item(18, 14).Value = 0.354
$ws.Cells.Item(18, 18).Value = 12.9
$ws.Cells.Item(18, 19).Value = 31
$ws.Cells.Item(18, 26).Value = 19
$ws.Cells.Item(18, 28).Value = 98.59999999999999
$ws.Cells.Item(18, 29).Value = -1.6
$ws.Cells.Item(18, 30).Value = 26
$ws.Cells.Item(18, 31).Value = 18
$ws.Cells.Item(18, 32).Value = 16
$ws.Cells.Item(18, 33).Value = 18
$ws.Cells.Item(18, 34).Value = 15
$ws.Cells.Item(18, 35).Value = 7
$ws.Cells.Item(18, 37).Value = 26
$ws.Cells.Item(19, 30).Value = 26
$ws.Cells.Item(19, 43).Value = 24
$ws.Cells.Item(20, 30).Value = 3
$ws.Cells.Item(21, 30).Value = 26
$ws.Cells.Item(21, 33).Value = 7
$ws.Cells.Item(22, 30).Value = 3
$ws.Cells.Item(22, 42).Value = 1
$ws.Cells.Item(23, 30).Value = 3
$ws.Cells.Item(24, 30).Value = 14
$ws.Cells.Item(25, 4).Value = 72
$ws.Cells.Item(25, 6).Value = 49
$ws.Cells.Item(25, 7).Value = 0.319
$ws.Cells.Item(25, 9).Value = 37
$ws.Cells.Item(25, 10).Value = 84.2
$ws.Cells.Item(25, 12).Value = 5.6
$ws.Cells.Item(25, 14).Value = 0.321
$ws.Cells.Item(25, 15).Value = 14.6
$ws.Cells.Item(25, 17).Value = 0.735
$ws.Cells.Item(25, 19).Value = 30
$ws.Cells.Item(25, 21).Value = 22.2
$ws.Cells.Item(25, 23).Value = 8
$ws.Cells.Item(25, 26).Value = 20.7
$ws.Cells.Item(25, 27).Value = 18.3
$ws.Cells.Item(25, 28).Value = 94.3
$ws.Cells.Item(25, 29).Value = -6.4
$ws.Cells.Item(25, 30).Value = 3
$ws.Cells.Item(25, 34).Value = 19
$ws.Cells.Item(25, 43).Value = 22
$ws.Cells.Item(26, 30).Value = 14
$ws.Cells.Item(27, 4).Value = 72
$ws.Cells.Item(27, 5).Value = 26
$ws.Cells.Item(27, 7).Value = 0.361
$ws.Cells.Item(27, 9).Value = 37.5
$ws.Cells.Item(27, 10).Value = 83.8
$ws.Cells.Item(27, 11).Value = 0.447
$ws.Cells.Item(27, 15).Value = 17.4
$ws.Cells.Item(27, 17).Value = 0.767
$ws.Cells.Item(27, 19).Value = 28.9
$ws.Cells.Item(27, 20).Value = 40.3
$ws.Cells.Item(27, 28).Value = 99.8
$ws.Cells.Item(27, 29).Value = -5
$ws.Cells.Item(27, 30).Value = 3
$ws.Cells.Item(27, 33).Value = 23
$ws.Cells.Item(27, 34).Value = 12
$ws.Cells.Item(27, 39).Value = 12
$ws.Cells.Item(27, 41).Value = 10
$ws.Cells.Item(27, 45).Value = 29
$ws.Cells.Item(27, 50).Value = 25
$ws.Cells.Item(27, 55).Value = 27
$ws.Cells.Item(28, 30).Value = 14
$ws.Cells.Item(29, 30).Value = 14
$ws.Cells.Item(29, 31).Value = 21
$ws.Cells.Item(29, 33).Value = 21
$ws.Cells.Item(29, 41).Value = 9
$ws.Cells.Item(30, 30).Value = 3
$ws.Cells.Item(30, 32).Value = 16
$ws.Cells.Item(30, 48).Value = 22
$ws.Cells.Item(31, 30).Value = 14
$ws.Cells.Item(31, 31).Value = 21
$ws.Cells.Item(31, 33).Value = 21
$ws.Cells.Item(31, 37).Value = 27
$ws.Cells.Item(31, 43).Value = 23

# --- Date column (BF) format fix: 'M-DD-YYYY-YY' -> 'YYYY-MM-DD' ---
# Force Text format first so Excel does not auto-convert the string into a date serial number.
$dateRange = $ws.Range($ws.Cells.Item(2, 58), $ws.Cells.Item(31, 58))
$dateRange.NumberFormat = "@"

$ws.Cells.Item(2, 58).Value = "2013-03-28"
$ws.Cells.Item(3, 58).Value = "2013-03-28"
$ws.Cells.Item(4, 58).Value = "2013-03-28"
$ws.Cells.Item(5, 58).Value = "2013-03-28"
$ws.Cells.Item(6, 58).Value = "2013-03-28"
$ws.Cells.Item(7, 58).Value = "2013-03-28"
$ws.Cells.Item(8, 58).Value = "2013-03-28"
$ws.Cells.Item(9, 58).Value = "2013-03-28"
$ws.Cells.Item(10, 58).Value = "2013-03-28"
$ws.Cells.Item(11, 58).Value = "2013-03-28"
$ws.Cells.Item(12, 58).Value = "2013-03-28"
$ws.Cells.Item(13, 58).Value = "2013-03-28"
$ws.Cells.Item(14, 58).Value = "2013-03-28"
$ws.Cells.Item(15, 58).Value = "2013-03-28"
$ws.Cells.Item(16, 58).Value = "2013-03-28"
$ws.Cells.Item(17, 58).Value = "2013-03-28"
$ws.Cells.Item(18, 58).Value = "2013-03-28"
$ws.Cells.Item(19, 58).Value = "2013-03-28"
$ws.Cells.Item(20, 58).Value = "2013-03-28"
$ws.Cells.Item(21, 58).Value = "2013-03-28"
$ws.Cells.Item(22, 58).Value = "2013-03-28"
$ws.Cells.Item(23, 58).Value = "2013-03-28"
$ws.Cells.Item(24, 58).Value = "2013-03-28"
$ws.Cells.Item(25, 58).Value = "2013-03-28"
$ws.Cells.Item(26, 58).Value = "2013-03-28"
$ws.Cells.Item(27, 58).Value = "2013-03-28"
$ws.Cells.Item(28, 58).Value = "2013-03-28"
$ws.Cells.Item(29, 58).Value = "2013-03-28"
$ws.Cells.Item(30, 58).Value = "2013-03-28"
$ws.Cells.Item(31, 58).Value = "2013-03-28"

# Clear the temporary number-format override so the cells end up with no explicit style,
# matching a plain text cell (keeps the value as text, not a date serial).
$dateRange.Style = "Normal"

